$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3332190829615296
$ws.Cells.Item(2, 3).Value = -0.8011660996346401
$ws.Cells.Item(2, 4).Value = 0.5787938608431268
$ws.Cells.Item(2, 5).Value = 0.1949980559329599
$ws.Cells.Item(2, 6).Value = 0.8798839467072684
$ws.Cells.Item(2, 7).Value = 0.1878585623255973
$ws.Cells.Item(2, 8).Value = 0.4273069063282261
$ws.Cells.Item(2, 9).Value = 0.4640348830873968
$ws.Cells.Item(2, 10).Value = 1.290528957474712
$ws.Cells.Item(2, 11).Value = 0.634730747353147

$ws.Cells.Item(3, 2).Value = 1.379959960477767
$ws.Cells.Item(3, 3).Value = 0.9961641555676001
$ws.Cells.Item(3, 4).Value = 1.681050046341908
$ws.Cells.Item(3, 5).Value = 0.9890246619602374
$ws.Cells.Item(3, 6).Value = 1.228473005962866
$ws.Cells.Item(3, 7).Value = 1.265200982722037
$ws.Cells.Item(3, 8).Value = 2.091695057109352
$ws.Cells.Item(3, 9).Value = 1.435896846987787
$ws.Cells.Item(3, 10).Value = 0.7777466072621664
$ws.Cells.Item(3, 11).Value = 1.176861555113949

$ws.Cells.Item(4, 2).Value = 0.6848858907743085
$ws.Cells.Item(4, 3).Value = -0.007139493607362657
$ws.Cells.Item(4, 4).Value = 0.2323088503952662
$ws.Cells.Item(4, 5).Value = 0.2690368271544369
$ws.Cells.Item(4, 6).Value = 1.095530901541752
$ws.Cells.Item(4, 7).Value = 0.4397326914201872
$ws.Cells.Item(4, 8).Value = -0.2184175483054337
$ws.Cells.Item(4, 9).Value = 0.1806973995463494
$ws.Cells.Item(4, 10).Value = 0.119260857424772
$ws.Cells.Item(4, 11).Value = -0.3149050776923316

$ws.Cells.Item(5, 2).Value = 0.2394483440026288
$ws.Cells.Item(5, 3).Value = 0.2761763207617995
$ws.Cells.Item(5, 4).Value = 1.102670395149115
$ws.Cells.Item(5, 5).Value = 0.4468721850275498
$ws.Cells.Item(5, 6).Value = -0.211278054698071
$ws.Cells.Item(5, 7).Value = 0.187836893153712
$ws.Cells.Item(5, 8).Value = 0.1264003510321347
$ws.Cells.Item(5, 9).Value = -0.3077655840849689
$ws.Cells.Item(5, 10).Value = 0.671035807994108
$ws.Cells.Item(5, 11).Value = 0.4439199035461818

$ws.Cells.Item(6, 2).Value = 0.8264940743873155
$ws.Cells.Item(6, 3).Value = 0.1706958642657503
$ws.Cells.Item(6, 4).Value = -0.4874543754598706
$ws.Cells.Item(6, 5).Value = -0.0883394276080875
$ws.Cells.Item(6, 6).Value = -0.1497759697296649
$ws.Cells.Item(6, 7).Value = -0.5839419048467684
$ws.Cells.Item(6, 8).Value = 0.3948594872323085
$ws.Cells.Item(6, 9).Value = 0.1677435827843823
$ws.Cells.Item(6, 10).Value = -0.4657161466417062
$ws.Cells.Item(6, 11).Value = -0.0383191879861019

$ws.Cells.Item(7, 2).Value = -0.6581502397256208
$ws.Cells.Item(7, 3).Value = -0.2590352918738378
$ws.Cells.Item(7, 4).Value = -0.3204718339954152
$ws.Cells.Item(7, 5).Value = -0.7546377691125187
$ws.Cells.Item(7, 6).Value = 0.2241636229665582
$ws.Cells.Item(7, 7).Value = -0.00295228148136796
$ws.Cells.Item(7, 8).Value = -0.6364120109074565
$ws.Cells.Item(7, 9).Value = -0.2090150522518522
$ws.Cells.Item(7, 10).Value = -0.3753520069399991
$ws.Cells.Item(7, 11).Value = -0.2526409904991733

$ws.Cells.Item(8, 2).Value = -0.0614365421215774
$ws.Cells.Item(8, 3).Value = -0.4956024772386809
$ws.Cells.Item(8, 4).Value = 0.483198914840396
$ws.Cells.Item(8, 5).Value = 0.2560830103924698
$ws.Cells.Item(8, 6).Value = -0.3773767190336187
$ws.Cells.Item(8, 7).Value = 0.0500202396219856
$ws.Cells.Item(8, 8).Value = -0.1163167150661613
$ws.Cells.Item(8, 9).Value = 0.006394301374664513
$ws.Cells.Item(8, 10).Value = 0.05362803792980481
$ws.Cells.Item(8, 11).Value = -0.3429906935926468

$ws.Cells.Item(9, 2).Value = 0.9788013920790769
$ws.Cells.Item(9, 3).Value = 0.7516854876311507
$ws.Cells.Item(9, 4).Value = 0.1182257582050622
$ws.Cells.Item(9, 5).Value = 0.5456227168606665
$ws.Cells.Item(9, 6).Value = 0.3792857621725196
$ws.Cells.Item(9, 7).Value = 0.5019967786133455
$ws.Cells.Item(9, 8).Value = 0.5492305151684858
$ws.Cells.Item(9, 9).Value = 0.1526117836460341
$ws.Cells.Item(9, 10).Value = 0.4066638187697542
$ws.Cells.Item(9, 11).Value = 0.6267882086433268

$ws.Cells.Item(10, 2).Value = -0.6334597294260885
$ws.Cells.Item(10, 3).Value = -0.2060627707704842
$ws.Cells.Item(10, 4).Value = -0.3723997254586311
$ws.Cells.Item(10, 5).Value = -0.2496887090178053
$ws.Cells.Item(10, 6).Value = -0.202454972462665
$ws.Cells.Item(10, 7).Value = -0.5990737039851166
$ws.Cells.Item(10, 8).Value = -0.3450216688613965
$ws.Cells.Item(10, 9).Value = -0.1248972789878239
$ws.Cells.Item(10, 10).Value = -0.3606857275744921
$ws.Cells.Item(10, 11).Value = -0.5396584850452768

$ws.Cells.Item(11, 2).Value = -0.1663369546881469
$ws.Cells.Item(11, 3).Value = -0.04362593824732108
$ws.Cells.Item(11, 4).Value = 0.00360779830781921
$ws.Cells.Item(11, 5).Value = -0.3930109332146324
$ws.Cells.Item(11, 6).Value = -0.1389588980909123
$ws.Cells.Item(11, 7).Value = 0.0811654917826603
$ws.Cells.Item(11, 8).Value = -0.1546229568040079
$ws.Cells.Item(11, 9).Value = -0.3335957142747926
$ws.Cells.Item(11, 10).Value = -0.3441562445341136
$ws.Cells.Item(11, 11).Value = -0.5518753546922437

$ws.Cells.Item(12, 2).Value = 0.04723373655514029
$ws.Cells.Item(12, 3).Value = -0.3493849949673113
$ws.Cells.Item(12, 4).Value = -0.09533295984359125
$ws.Cells.Item(12, 5).Value = 0.1247914300299814
$ws.Cells.Item(12, 6).Value = -0.1109970185566868
$ws.Cells.Item(12, 7).Value = -0.2899697760274715
$ws.Cells.Item(12, 8).Value = -0.3005303062867926
$ws.Cells.Item(12, 9).Value = -0.5082494164449226
$ws.Cells.Item(12, 10).Value = -0.2910405169700179
$ws.Cells.Item(12, 11).Value = 0.0186456522880436

$ws.Cells.Item(13, 2).Value = 0.2540520351237201
$ws.Cells.Item(13, 3).Value = 0.4741764249972927
$ws.Cells.Item(13, 4).Value = 0.2383879764106245
$ws.Cells.Item(13, 5).Value = 0.05941521893983981
$ws.Cells.Item(13, 6).Value = 0.04885468868051879
$ws.Cells.Item(13, 7).Value = -0.1588644214776113
$ws.Cells.Item(13, 8).Value = 0.05834447799729348
$ws.Cells.Item(13, 9).Value = 0.3680306472553549
$ws.Cells.Item(13, 10).Value = 0.2137265134233509
$ws.Cells.Item(13, 11).Value = 0.8641529346425341

$ws.Cells.Item(14, 2).Value = -0.2357884485866682
$ws.Cells.Item(14, 3).Value = -0.4147612060574529
$ws.Cells.Item(14, 4).Value = -0.4253217363167739
$ws.Cells.Item(14, 5).Value = -0.633040846474904
$ws.Cells.Item(14, 6).Value = -0.4158319469999993
$ws.Cells.Item(14, 7).Value = -0.1061457777419378
$ws.Cells.Item(14, 8).Value = -0.2604499115739418
$ws.Cells.Item(14, 9).Value = 0.3899765096452414
$ws.Cells.Item(14, 10).Value = 0.1738114610270322
$ws.Cells.Item(14, 11).Value = -0.2332176680079241

$ws.Cells.Item(15, 2).Value = -0.01056053025932102
$ws.Cells.Item(15, 3).Value = -0.2182796404174511
$ws.Cells.Item(15, 4).Value = -0.001070740942546333
$ws.Cells.Item(15, 5).Value = 0.3086154283155151
$ws.Cells.Item(15, 6).Value = 0.1543112944835111
$ws.Cells.Item(15, 7).Value = 0.8047377157026943
$ws.Cells.Item(15, 8).Value = 0.5885726670844852
$ws.Cells.Item(15, 9).Value = 0.1815435380495288
$ws.Cells.Item(15, 10).Value = 0.6805857552269486
$ws.Cells.Item(15, 11).Value = 0.5217368420714361

$ws.Cells.Item(16, 2).Value = 0.2172088994749047
$ws.Cells.Item(16, 3).Value = 0.5268950687329662
$ws.Cells.Item(16, 4).Value = 0.3725909349009622
$ws.Cells.Item(16, 5).Value = 1.023017356120145
$ws.Cells.Item(16, 6).Value = 0.8068523075019363
$ws.Cells.Item(16, 7).Value = 0.3998231784669799
$ws.Cells.Item(16, 8).Value = 0.8988653956443997
$ws.Cells.Item(16, 9).Value = 0.7400164824888872
$ws.Cells.Item(16, 10).Value = 0.7961803034591312
$ws.Cells.Item(16, 11).Value = 3.045983143070524

$ws.Cells.Item(17, 2).Value = 0.3096861692580615
$ws.Cells.Item(17, 3).Value = 0.1553820354260574
$ws.Cells.Item(17, 4).Value = 0.8058084566452406
$ws.Cells.Item(17, 5).Value = 0.5896434080270315
$ws.Cells.Item(17, 6).Value = 0.1826142789920752
$ws.Cells.Item(17, 7).Value = 0.681656496169495
$ws.Cells.Item(17, 8).Value = 0.5228075830139824
$ws.Cells.Item(17, 9).Value = 0.5789714039842264
$ws.Cells.Item(17, 10).Value = 2.828774243595618
$ws.Cells.Item(17, 11).Value = 10.29792457445291

$ws.Cells.Item(18, 2).Value = -0.154304133832004
$ws.Cells.Item(18, 3).Value = 0.4961222873871792
$ws.Cells.Item(18, 4).Value = 0.2799572387689701
$ws.Cells.Item(18, 5).Value = -0.1270718902659863
$ws.Cells.Item(18, 6).Value = 0.3719703269114335
$ws.Cells.Item(18, 7).Value = 0.2131214137559209
$ws.Cells.Item(18, 8).Value = 0.269285234726165
$ws.Cells.Item(18, 9).Value = 2.519088074337557
$ws.Cells.Item(18, 10).Value = 9.988238405194851
$ws.Cells.Item(18, 11).Value = -8.183023394942618

$ws.Cells.Item(19, 2).Value = 0.6504264212191833
$ws.Cells.Item(19, 3).Value = 0.4342613726009741
$ws.Cells.Item(19, 4).Value = 0.02723224356601772
$ws.Cells.Item(19, 5).Value = 0.5262744607434375
$ws.Cells.Item(19, 6).Value = 0.3674255475879249
$ws.Cells.Item(19, 7).Value = 0.423589368558169
$ws.Cells.Item(19, 8).Value = 2.673392208169561
$ws.Cells.Item(19, 9).Value = 10.14254253902685
$ws.Cells.Item(19, 10).Value = -8.028719261110615
$ws.Cells.Item(19, 11).Value = 0.1109522858685349

$ws.Cells.Item(20, 2).Value = -0.2161650486182091
$ws.Cells.Item(20, 3).Value = -0.6231941776531655
$ws.Cells.Item(20, 4).Value = -0.1241519604757457
$ws.Cells.Item(20, 5).Value = -0.2830008736312583
$ws.Cells.Item(20, 6).Value = -0.2268370526610142
$ws.Cells.Item(20, 7).Value = 2.022965786950378
$ws.Cells.Item(20, 8).Value = 9.49211611780767
$ws.Cells.Item(20, 9).Value = -8.679145682329798
$ws.Cells.Item(20, 10).Value = -0.5394741353506483
$ws.Cells.Item(20, 11).Value = 1.608096336448744

$ws.Cells.Item(21, 2).Value = -0.4070291290349564
$ws.Cells.Item(21, 3).Value = 0.09201308814246346
$ws.Cells.Item(21, 4).Value = -0.06683582501304909
$ws.Cells.Item(21, 5).Value = -0.01067200404280504
$ws.Cells.Item(21, 6).Value = 2.239130835568587
$ws.Cells.Item(21, 7).Value = 9.70828116642588
$ws.Cells.Item(21, 8).Value = -8.462980633711588
$ws.Cells.Item(21, 9).Value = -0.3233090867324392
$ws.Cells.Item(21, 10).Value = 1.824261385066953
$ws.Cells.Item(21, 11).Value = -1.615768601456377

$ws.Cells.Item(22, 2).Value = 0.4990422171774198
$ws.Cells.Item(22, 3).Value = 0.3401933040219072
$ws.Cells.Item(22, 4).Value = 0.3963571249921513
$ws.Cells.Item(22, 5).Value = 2.646159964603544
$ws.Cells.Item(22, 6).Value = 10.11531029546084
$ws.Cells.Item(22, 7).Value = -8.055951504676631
$ws.Cells.Item(22, 8).Value = 0.08372004230251717
$ws.Cells.Item(22, 9).Value = 2.231290514101909
$ws.Cells.Item(22, 10).Value = -1.20873947242142
$ws.Cells.Item(22, 11).Value = -1.269913650835621

$ws.Cells.Item(23, 2).Value = -0.1588489131555126
$ws.Cells.Item(23, 3).Value = -0.1026850921852685
$ws.Cells.Item(23, 4).Value = 2.147117747426124
$ws.Cells.Item(23, 5).Value = 9.616268078283417
$ws.Cells.Item(23, 6).Value = -8.554993721854052
$ws.Cells.Item(23, 7).Value = -0.4153221748749026
$ws.Cells.Item(23, 8).Value = 1.73224829692449
$ws.Cells.Item(23, 9).Value = -1.70778168959884
$ws.Cells.Item(23, 10).Value = -1.768955868013041
$ws.Cells.Item(23, 11).Value = 0.2725000691432133

$ws.Cells.Item(24, 2).Value = 0.05616382097024405
$ws.Cells.Item(24, 3).Value = 2.305966660581636
$ws.Cells.Item(24, 4).Value = 9.77511699143893
$ws.Cells.Item(24, 5).Value = -8.396144808698539
$ws.Cells.Item(24, 6).Value = -0.2564732617193901
$ws.Cells.Item(24, 7).Value = 1.891097210080002
$ws.Cells.Item(24, 8).Value = -1.548932776443328
$ws.Cells.Item(24, 9).Value = -1.610106954857529
$ws.Cells.Item(24, 10).Value = 0.4313489822987259
$ws.Cells.Item(24, 11).Value = -0.1512455547349151

